# Applies the "math_L-curve" perturbation to the optimization_parameters sheet:
#  - removes the "Deletion" row
#  - renames the "Model" label to "production_function"
#  - adds a new "L_curve" parameter row (value 1) right after production_function
#  - trims the redundant C1:F1 header filler cells on row 1
#  - updates the sheet's active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Remove the now-unused repeated "value" header cells in C1:F1 (only A1/B1 remain).
$ws.Range("C1:F1").ClearContents()

# Remove the "Deletion" row (row 16: "Deletion" / 0 / 3).
$ws.Rows.Item(16).Delete()

# Rename "Model" (row 8, column A) to "production_function".
$ws.Range("A8").Value = "production_function"

# Insert a new row below "production_function" for the "L_curve" parameter.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 1

# Update the sheet's active selection to match the reformatted layout.
$ws.Range("C1:G5").Select()
